$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column D (Description shifts from D to E)
$ws.Range("D1").EntireColumn.Insert()

# Match column D's width to column C so they merge into one col span
$ws.Range("D1").ColumnWidth = 16.0221354166667

# New header for column D
$ws.Range("D2").Value = "Status"

# Status values for specific rows
$ws.Range("D4").Value = "Done"
$ws.Range("D6").Value = "Done"
$ws.Range("D10").Value = "In progress"

# Update selection to match the target state
$ws.Range("D5").Select()
